# This script updates the "想去人数" (number of people wanting to go) column (F)
# for a handful of exhibition rows on both the "展览" sheet and the aggregated
# "全部类型" sheet, matching the data refresh recorded in the commit
# "Update gh-pages to output generated at 456a3b4".

$wb = $excel.ActiveWorkbook

# Column F updates for sheet "展览" (row -> new value)
$exhibitionUpdates = @{
    9  = 722
    10 = 2189
    12 = 1658
    13 = 2783
    14 = 146
    15 = 4159
    17 = 174
    19 = 528
    24 = 285
    25 = 4046
    27 = 3504
    28 = 1111
    34 = 439
}

$wsExhibition = $wb.Worksheets.Item("展览")
foreach ($row in $exhibitionUpdates.Keys) {
    $wsExhibition.Range("F$row").Value = $exhibitionUpdates[$row]
}

# Column F updates for sheet "全部类型" (row -> new value)
$allTypesUpdates = @{
    11 = 722
    12 = 2189
    14 = 1658
    16 = 2783
    17 = 146
    18 = 4159
    20 = 174
    22 = 528
    28 = 285
    29 = 4046
    31 = 3504
    32 = 1111
    38 = 439
}

$wsAllTypes = $wb.Worksheets.Item("全部类型")
foreach ($row in $allTypesUpdates.Keys) {
    $wsAllTypes.Range("F$row").Value = $allTypesUpdates[$row]
}

$wb.Save()
